$d = $word.ActiveDocument

function Set-ParagraphXml($range, $innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $innerXml +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- Paragraph 1: "Weather In California" -> split "In" out with grammar-check marks ---
$p1 = '<w:body><w:p>' +
      '<w:r><w:t xml:space="preserve">Weather </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>In</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:t xml:space="preserve"> California</w:t></w:r>' +
      '</w:p></w:body>'
Set-ParagraphXml $d.Paragraphs.Item(1).Range $p1

# --- Paragraph 3: merge the bookmark-split runs into a single run, drop the _GoBack bookmark here ---
$p3 = '<w:body><w:p>' +
      '<w:r><w:t>Yesterday' + [char]0x2019 + 's Weather: As my dad put it, ' + [char]0x201C + 'Perfect day to go to the beach and play volley ball with bikini models!' + [char]0x201D + ' ' + [char]0x2026 + '</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> ok dad thanks -_-</w:t></w:r>' +
      '</w:p></w:body>'
Set-ParagraphXml $d.Paragraphs.Item(3).Range $p3

# --- Paragraph 5: mark "gonna" / "wanna" as spelling errors ---
$p5 = '<w:body><w:p>' +
      '<w:r><w:t xml:space="preserve">Today' + [char]0x2019 + 's Weather: As my dad would say, ' + [char]0x201C + 'When are you </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>gonna</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> visit? Aren' + [char]0x2019 + 't you tired of the snow?' + [char]0x201D + ' Shut up dad I don' + [char]0x2019 + 't </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>wanna</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> hear how I' + [char]0x2019 + 'm missing out because it gets 30 below here.</w:t></w:r>' +
      '</w:p></w:body>'
Set-ParagraphXml $d.Paragraphs.Item(5).Range $p5

# --- Append the new "Florida" section at the end of the document ---
$newParas = '<w:body>' +
      '<w:p/>' +
      '<w:p>' +
      '<w:r><w:t xml:space="preserve">Weather </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>In</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:t xml:space="preserve"> Florida</w:t></w:r>' +
      '</w:p>' +
      '<w:p/>' +
      '<w:p>' +
      '<w:r><w:t xml:space="preserve">Yesterday' + [char]0x2019 + 's weather: Not </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>gonna</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> lie it was super humid.</w:t></w:r>' +
      '</w:p>' +
      '<w:p>' +
      '<w:r><w:t>Today' + [char]0x2019 + 's weather: According to my brother' + [char]0x2026 + ' ' + [char]0x201C + 'Remember when you accidentally touched a car lighter?' + [char]0x201D + '</w:t></w:r>' +
      '</w:p>' +
      '<w:p>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
      '</w:p>' +
      '</w:body>'

$endRange = $d.Content
$endRange.Collapse(0)
Set-ParagraphXml $endRange $newParas
